$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 14 (Boston River vs Rampla Juniors) so that the former
# row 15 (Wanderers vs Penarol) shifts up to become the new row 14.
$ws.Rows.Item(14).Delete()

# Update individual odds cells that changed value in the remaining rows.

# Row 5 (Dep. Pasto vs Patriotas)
$ws.Range("G5").Value = 1.33
$ws.Range("I5").Value = 12
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("W5").Value = 5
$ws.Range("AC5").Value = 7.5
$ws.Range("AF5").Value = 126
$ws.Range("AH5").Value = 19
$ws.Range("AJ5").Value = 34
$ws.Range("AK5").Value = 151
$ws.Range("AM5").Value = 101
$ws.Range("AU5").Value = 12
$ws.Range("AY5").Value = 51

# Row 9 (Junior vs Dep. Cali)
$ws.Range("G9").Value = 1.42
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 8.5
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2.25
$ws.Range("L9").Value = 8
$ws.Range("U9").Value = 2.25
$ws.Range("V9").Value = 1.57
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 9
$ws.Range("AE9").Value = 23
$ws.Range("AH9").Value = 19
$ws.Range("AJ9").Value = 26
$ws.Range("AK9").Value = 101
$ws.Range("AL9").Value = 67
$ws.Range("AN9").Value = 3.2
$ws.Range("AO9").Value = 7
$ws.Range("AQ9").Value = 21
$ws.Range("AZ9").Value = 201
$ws.Range("BA9").Value = 251

# Row 12 (Celaya vs Atlante)
$ws.Range("N12").Value = 4.95

# Row 13 (Leones Negros vs Tapatio)
$ws.Range("G13").Value = 1.7
$ws.Range("I13").Value = 5.1
$ws.Range("J13").Value = 2.27
$ws.Range("L13").Value = 5.3
$ws.Range("N13").Value = 9.15
$ws.Range("R13").Value = 1.83
$ws.Range("X13").Value = 8.25
$ws.Range("Y13").Value = 7.8
$ws.Range("Z13").Value = 14
$ws.Range("AA13").Value = 13
$ws.Range("AB13").Value = 24
$ws.Range("AC13").Value = 9.25
$ws.Range("AE13").Value = 14.5
$ws.Range("AF13").Value = 70
$ws.Range("AH13").Value = 13
$ws.Range("AI13").Value = 30
$ws.Range("AJ13").Value = 16
$ws.Range("AK13").Value = 100
$ws.Range("AL13").Value = 55
$ws.Range("AM13").Value = 55
$ws.Range("AN13").Value = 3.5
$ws.Range("AO13").Value = 8.5
$ws.Range("AP13").Value = 17
$ws.Range("AQ13").Value = 29
$ws.Range("AS13").Value = 200
$ws.Range("AU13").Value = 7.2
$ws.Range("AW13").Value = 6.7
$ws.Range("AX13").Value = 32
$ws.Range("AY13").Value = 35
$ws.Range("AZ13").Value = 200
$ws.Range("BA13").Value = 250
$ws.Range("BB13").Value = 500
